{"js": "// Replace each unique problem/date string in the document body with its\n// updated value, matching the OOXML diff (date header + 25 multiplication\n// problems across the practice table).\nconst replacements = [\n  [\"2024-11-08 Friday\", \"2024-11-09 Saturday\"],\n  [\"60\u00d760=\", \"97\u00d761=\"],\n  [\"88\u00d723=\", \"79\u00d742=\"],\n  [\"76\u00d729=\", \"35\u00d727=\"],\n  [\"35\u00d733=\", \"97\u00d729=\"],\n  [\"71\u00d792=\", \"80\u00d733=\"],\n  [\"84\u00d795=\", \"94\u00d763=\"],\n  [\"27\u00d723=\", \"70\u00d725=\"],\n  [\"74\u00d726=\", \"57\u00d772=\"],\n  [\"17\u00d792=\", \"34\u00d790=\"],\n  [\"87\u00d747=\", \"96\u00d768=\"],\n  [\"88\u00d717=\", \"73\u00d754=\"],\n  [\"41\u00d752=\", \"22\u00d751=\"],\n  [\"80\u00d778=\", \"51\u00d742=\"],\n  [\"67\u00d723=\", \"62\u00d781=\"],\n  [\"73\u00d788=\", \"83\u00d730=\"],\n  [\"38\u00d795=\", \"39\u00d746=\"],\n  [\"32\u00d739=\", \"92\u00d723=\"],\n  [\"25\u00d744=\", \"68\u00d779=\"],\n  [\"23\u00d719=\", \"88\u00d754=\"],\n  [\"94\u00d751=\", \"26\u00d731=\"],\n  [\"46\u00d774=\", \"64\u00d718=\"],\n  [\"61\u00d764=\", \"33\u00d765=\"],\n  [\"28\u00d763=\", \"53\u00d713=\"],\n  [\"38\u00d759=\", \"80\u00d725=\"],\n  [\"28\u00d717=\", \"86\u00d799=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the practice sheet: the date header and all 25 two-digit by\n# two-digit multiplication problems in the table get new values, matching\n# the OOXML diff. Each old problem string is unique in the document, so a\n# simple Find/Replace (wdReplaceAll) per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-08 Friday\", \"2024-11-09 Saturday\"),\n    @(\"60\u00d760=\", \"97\u00d761=\"),\n    @(\"88\u00d723=\", \"79\u00d742=\"),\n    @(\"76\u00d729=\", \"35\u00d727=\"),\n    @(\"35\u00d733=\", \"97\u00d729=\"),\n    @(\"71\u00d792=\", \"80\u00d733=\"),\n    @(\"84\u00d795=\", \"94\u00d763=\"),\n    @(\"27\u00d723=\", \"70\u00d725=\"),\n    @(\"74\u00d726=\", \"57\u00d772=\"),\n    @(\"17\u00d792=\", \"34\u00d790=\"),\n    @(\"87\u00d747=\", \"96\u00d768=\"),\n    @(\"88\u00d717=\", \"73\u00d754=\"),\n    @(\"41\u00d752=\", \"22\u00d751=\"),\n    @(\"80\u00d778=\", \"51\u00d742=\"),\n    @(\"67\u00d723=\", \"62\u00d781=\"),\n    @(\"73\u00d788=\", \"83\u00d730=\"),\n    @(\"38\u00d795=\", \"39\u00d746=\"),\n    @(\"32\u00d739=\", \"92\u00d723=\"),\n    @(\"25\u00d744=\", \"68\u00d779=\"),\n    @(\"23\u00d719=\", \"88\u00d754=\"),\n    @(\"94\u00d751=\", \"26\u00d731=\"),\n    @(\"46\u00d774=\", \"64\u00d718=\"),\n    @(\"61\u00d764=\", \"33\u00d765=\"),\n    @(\"28\u00d763=\", \"53\u00d713=\"),\n    @(\"38\u00d759=\", \"80\u00d725=\"),\n    @(\"28\u00d717=\", \"86\u00d799=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $pair[0],   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap = wdFindContinue\n        $false,     # Format\n        $pair[1],   # ReplaceWith\n        2           # Replace = wdReplaceAll\n    ) | Out-Null\n}\n"}
